$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The text-like columns (A-D) auto-convert to dates/numbers when assigned
# a bare string that looks like one, so force text format before writing,
# then clear the explicit formatting back off (keeping the cell's type as
# text) so no extra style index is left behind on the new row.
$textRange = $ws.Range("A36:D36")
$textRange.NumberFormat = "@"

$ws.Range("A36").Value = "2023-06-09"
$ws.Range("B36").Value = "18:41:57"
$ws.Range("C36").Value = "Friday"
$ws.Range("D36").Value = "23"

$textRange.ClearFormats()

$ws.Range("E36").Value = 120860
$ws.Range("F36").Value = 134500
$ws.Range("G36").Value = 160764
$ws.Range("H36").Value = 132025
$ws.Range("I36").Value = 175890
$ws.Range("J36").Value = 113712
$ws.Range("K36").Value = 201670
$ws.Range("L36").Value = 221809
$ws.Range("M36").Value = 173245
$ws.Range("N36").Value = 118792
$ws.Range("O36").Value = 38669
$ws.Range("P36").Value = 34331
$ws.Range("Q36").Value = 51075
$ws.Range("R36").Value = -1
$ws.Range("S36").Value = 36877
$ws.Range("T36").Value = -1
